# Generate Report for Handoff
# Updates the localization-status report: the a50ba3c3 and b6657c0d source
# files are now "Ready for handoff" (new handoff cycle kicked off), with
# fresh handoff timestamps and a warning that the existing handback is
# stale relative to the newest source commit.

$wb = $excel.ActiveWorkbook

$statusOld = "Handed back: in sync with en-US"
$statusNew = "Ready for handoff"

$errA50 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/29d08affcf480cddcf396083ad8eb287b2c2be68/e2e/a50ba3c3-a210-4225-8a61-f5b79a37f6af.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d507287d3127133b504890e46d4e23b1675716b1/e2e/a50ba3c3-a210-4225-8a61-f5b79a37f6af.md."
$errB66 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/29d08affcf480cddcf396083ad8eb287b2c2be68/e2e/b6657c0d-4cd5-4065-82f2-45cedfeb434f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d507287d3127133b504890e46d4e23b1675716b1/e2e/b6657c0d-4cd5-4065-82f2-45cedfeb434f.md."

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E4").Value = $statusNew
$wsOverview.Range("F4").Value = $statusNew
$wsOverview.Range("G4").Value = "2016-09-06 10:34:32"

$wsOverview.Range("E5").Value = $statusNew
$wsOverview.Range("F5").Value = $statusNew
$wsOverview.Range("G5").Value = "2016-09-06 10:34:32"

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C4").Value = $statusNew
$wsZhCn.Range("H4").Value = "2016-09-06 10:34:28"
$wsZhCn.Range("P4").Value = $errA50

$wsZhCn.Range("C5").Value = $statusNew
$wsZhCn.Range("H5").Value = "2016-09-06 10:34:28"
$wsZhCn.Range("P5").Value = $errB66

$wsZhCn.Columns.Item(16).ColumnWidth = 39.16666666666667

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C4").Value = $statusNew
$wsDeDe.Range("H4").Value = "2016-09-06 10:34:32"
$wsDeDe.Range("P4").Value = $errA50

$wsDeDe.Range("C5").Value = $statusNew
$wsDeDe.Range("H5").Value = "2016-09-06 10:34:32"
$wsDeDe.Range("P5").Value = $errB66

$wsDeDe.Columns.Item(16).ColumnWidth = 39.16666666666667
